$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 734, shifting existing rows 734:775 down to 735:776
$ws.Rows.Item(734).Insert()

# Populate the newly inserted row 734 with the new data point.
# Column A holds a date-looking string that must stay literal text (not be
# auto-converted into a real date serial number), so force Text format
# before assigning the value, then clear the format again so the cell ends
# up with the same (default) style as its neighbours.
$ws.Cells.Item(734, 1).NumberFormat = "@"
$ws.Cells.Item(734, 1).Value = "2026/02/01"
$ws.Cells.Item(734, 1).ClearFormats()

$ws.Cells.Item(734, 2).Value = "日"
$ws.Cells.Item(734, 3).Value = 7
$ws.Cells.Item(734, 4).Value = 23
